$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Copy row formatting/styles for the new rows BEFORE changing the
#        source rows' own content, so the styles (B13:D13 header style,
#        B14:I14 data style) are captured correctly. ---

# New section-header rows (15 and 17) get row 13's header style (B:D only)
$ws.Range("B13:D13").Copy($ws.Range("B15:D15"))
$ws.Range("B13:D13").Copy($ws.Range("B17:D17"))

# New data rows (16 and 18) get row 14's full data-row style (B:I)
$ws.Range("B14:I14").Copy($ws.Range("B16:I16"))
$ws.Range("B14:I14").Copy($ws.Range("B18:I18"))

# --- 2) Update existing data rows 10-12: Ejecución (G) and Avance (H) ---
$ws.Range("G10").Value = 0.66
$ws.Range("H10").Value = 1
$ws.Range("G11").Value = 153.07
$ws.Range("H11").Value = 1
$ws.Range("G12").Value = 5.09
$ws.Range("H12").Value = 1

# --- 3) Update row 13 (category header) to MR200 / Limpieza de obras de arte ---
$ws.Range("C13").Value = "MR200"
$ws.Range("D13").Value = "Limpieza de obras de arte"

# --- 4) Update row 14 (item 4) to MR203 / Limpieza de badén ---
$ws.Range("C14").Value = "MR203"
$ws.Range("D14").Value = "Limpieza de badén"
$ws.Range("E14").Value = "m2"
$ws.Range("F14").Value = 22.224
$ws.Range("G14").Value = 22.22
$ws.Range("H14").Value = 1

# --- 5) Row 15 (category header) = MR300 / Control de vegetación ---
$ws.Range("C15").Value = "MR300"
$ws.Range("D15").Value = "Control de vegetación"

# --- 6) Row 16 (item 5) = MR301 / Roce y limpieza ---
$ws.Range("B16").Value = 5
$ws.Range("C16").Value = "MR301"
$ws.Range("D16").Value = "Roce y limpieza"
$ws.Range("E16").Value = "m2"
$ws.Range("F16").Value = 1951.915
$ws.Range("G16").Value = 1951.52
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = "0+000 - 11+920"

# --- 7) Row 17 (category header) = MR400 / Seguridad vial ---
$ws.Range("C17").Value = "MR400"
$ws.Range("D17").Value = "Seguridad vial"

# --- 8) Row 18 (item 6) = MR401 / Conservación de señales ---
$ws.Range("B18").Value = 6
$ws.Range("C18").Value = "MR401"
$ws.Range("D18").Value = "Conservación de señales"
$ws.Range("E18").Value = "unidad"
$ws.Range("F18").Value = 5.323333333333333
$ws.Range("G18").Value = 5.32
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = "0+000 - 11+920"

# --- 9) Extend the conditional formatting range to cover the new rows ---
$fc = $ws.Range("B7:I14").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("B7:I18"))
